$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows covering 2021-09-02 through 2021-09-09 (serials 44441-44448)
# Columns: A = date (serial), B = nuovi pos., C = somma mobile 7gg., D = somma mobile 7gg. per 100mila abitanti
$rows = @(
  @(367, 44441, 2, 12, 129.2546316242999),
  @(368, 44442, 3, 14, 150.7970702283498),
  @(369, 44443, 3, 15, 161.5682895303748),
  @(370, 44444, 4, 14, 150.7970702283498),
  @(371, 44445, 2, 15, 161.5682895303748),
  @(372, 44446, 1, 16, 172.3395088323998),
  @(373, 44447, 0, 15, 161.5682895303748),
  @(374, 44448, 2, 15, 161.5682895303748)
)

foreach ($r in $rows) {
  $rowNum = $r[0]

  # Carry the date-column style (numFmt, border, font, alignment) down from
  # the previous row so the new cells match the existing column A formatting.
  $ws.Range("A366").Copy()
  $ws.Range("A" + $rowNum).PasteSpecial(-4122)

  $ws.Range("A" + $rowNum).Value = $r[1]
  $ws.Range("B" + $rowNum).Value = $r[2]
  $ws.Range("C" + $rowNum).Value = $r[3]
  $ws.Range("D" + $rowNum).Value = $r[4]
}

$excel.CutCopyMode = $false
